# ToR - original web client
# Slide 2: update title + content placeholder with the "True North Surf Club" brief.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Title -------------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "True North Surf Club"

# --- Content placeholder -------------------------------------------------
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# Paragraph 1: replace "- Information about client." in place (keeps the
# paragraph mark / subsequent paragraphs untouched).
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "True North Surf Club is an independent Surf Club that is affiliated with Surfing England, and is based in North Tyneside."

# Paragraph 2: replace "- Their initial 'pitch' ..." in place.
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "They have requested the creation of a website for the surf club, which includes the following features:"

# Insert four brand-new bullet paragraphs right after paragraph 2, before
# "The team consists of 5 members...". Leading "`r" starts a fresh
# paragraph without disturbing paragraph 2's own text/mark.
$newBlock = "- a gallery page to showcase their images.`r- a contact form to provide a means of contacting the club.`r- a membership sign up page, allowing potential new members to sign up for membership.`r- an e-commerce page so that branded merchandise can be purchased from the club via the website."
$para2.InsertAfter("`r" + $newBlock) | Out-Null

# Re-fetch the freshly created paragraphs (3..6) and bold the requested
# sub-phrases within each one.
$p3 = $tr.Paragraphs(3, 1)
$tr.Characters($p3.Start + 4, 12).Font.Bold = $true   # "gallery page"

$p4 = $tr.Paragraphs(4, 1)
$tr.Characters($p4.Start + 4, 13).Font.Bold = $true   # "contact form "

$p5 = $tr.Paragraphs(5, 1)
$tr.Characters($p5.Start + 4, 23).Font.Bold = $true   # "membership sign up page"

$p6 = $tr.Paragraphs(6, 1)
$tr.Characters($p6.Start + 5, 16).Font.Bold = $true   # "e-commerce page "

# --- Shrink text to fit the placeholder (as PowerPoint does once the new
# paragraphs overflow it). ------------------------------------------------
$body.TextFrame.AutoSize = 2
